$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update leaderboard rows 1-4 with new names/scores
$ws.Range("A1").Value = "Jack"
$ws.Range("B1").Value = 9
$ws.Range("C1").Value = 5
$ws.Range("D1").Value = 1.8

$ws.Range("A2").Value = "Lance"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 1.333

$ws.Range("A3").Value = "Lance "
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1

$ws.Range("A4").Value = "Connor"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 3
